$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: keep values, only change project_count (F2) from 5 to 1
$ws.Range("F2").Value = 1

# Row 3: replace with SMART SENSING MIDDLEWARE data
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "SMART SENSING MIDDLEWARE"
$ws.Range("C3").Value = 85.02
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 85.02
$ws.Range("F3").Value = 1

# Row 4: replace with Website for the Literature Society of the college data
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Website for the Literature Society of the college"
$ws.Range("C4").Value = 81.62
$ws.Range("D4").Value = 0.9
$ws.Range("E4").Value = 73.45999999999999
$ws.Range("F4").Value = 2

# Row 5: replace with LLMGuard data
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "LLMGuard"
$ws.Range("C5").Value = 78.28
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 78.28
$ws.Range("F5").Value = 2

# Remove rows 6-13 (old leftover data) without corrupting sheet - use ClearContents
$ws.Range("A6:F13").ClearContents()
